# "bold text on diagram"
# Slide 3 contains the flowchart diamonds. Bold the leading "Hash similarity"
# portion of the first diamond's text (leaving the trailing "?" unbolded,
# which splits the run in two), and bold the full "Tokenized yet?" text in
# the other two diamonds.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(3)

# Diamond 6: "Hash similarity?" -> bold "Hash similarity", keep "?" unbolded
$sh = $s.Shapes.Item(10)
$tr = $sh.TextFrame.TextRange
$tr.Characters(1, 15).Font.Bold = 1

# Diamond 36: "Tokenized yet?" -> bold entire text
$sh = $s.Shapes.Item(15)
$sh.TextFrame.TextRange.Font.Bold = 1

# Diamond 55: "Tokenized yet?" -> bold entire text
$sh = $s.Shapes.Item(20)
$sh.TextFrame.TextRange.Font.Bold = 1
